$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report was regenerated: the old handback run referenced markdown
# files "1e084702-...md" and "47d8a4f5-...md"; the new run references
# "a901f226-...md" and "ffffa8a81730-...md" (with fresh xliff hashes and
# updated timestamps). Update every sheet (Overview, zh-cn, de-de) plus the
# hyperlink display text that mirrors those file names.
# ---------------------------------------------------------------------------

# ---------------- Overview sheet ----------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws1.Range("B2").Value = "e2e\a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws1.Range("G2").Value = "2016-08-12 15:16:58"

$ws1.Range("A3").Value = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws1.Range("B3").Value = "e2e\ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws1.Range("G3").Value = "2016-08-12 15:16:58"

foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.TextToDisplay -eq "e2e\1e084702-b98b-41f1-9618-c1d5eeb137cc.md") {
        $hl.TextToDisplay = "e2e\a901f226-5706-402d-ac11-32b1cc14ef39.md"
    } elseif ($hl.TextToDisplay -eq "e2e\47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md") {
        $hl.TextToDisplay = "e2e\ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
    }
}

# ---------------- zh-cn sheet ----------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws2.Range("G2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-12 15:16:51"
$ws2.Range("I2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws2.Range("J2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-12 15:17:19"

$ws2.Range("A3").Value = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws2.Range("G3").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-12 15:16:51"
$ws2.Range("I3").Value = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws2.Range("J3").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-12 15:17:19"

foreach ($hl in $ws2.Hyperlinks) {
    if ($hl.TextToDisplay -eq "1e084702-b98b-41f1-9618-c1d5eeb137cc.md") {
        $hl.TextToDisplay = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
    } elseif ($hl.TextToDisplay -eq "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md") {
        $hl.TextToDisplay = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
    }
}

# ---------------- de-de sheet ----------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws3.Range("G2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-12 15:16:58"
$ws3.Range("I2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
$ws3.Range("J2").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-12 15:17:28"

$ws3.Range("A3").Value = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws3.Range("G3").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-12 15:16:58"
$ws3.Range("I3").Value = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
$ws3.Range("J3").Value = "a901f226-5706-402d-ac11-32b1cc14ef39.3f72c0d8e83c84a169663c64289b64ee5536ecce.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-12 15:17:28"

foreach ($hl in $ws3.Hyperlinks) {
    if ($hl.TextToDisplay -eq "1e084702-b98b-41f1-9618-c1d5eeb137cc.md") {
        $hl.TextToDisplay = "a901f226-5706-402d-ac11-32b1cc14ef39.md"
    } elseif ($hl.TextToDisplay -eq "47d8a4f5-dabe-4ee9-91be-54338ba3fb95.md") {
        $hl.TextToDisplay = "ffffa8a81730-ccd7-4ce7-a539-c51bf381a860.md"
    }
}
